$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 6 new rows at row 3785, shifting existing rows 3785:3954 down to 3791:3960
$ws.Rows("3785:3790").Insert()

# Populate the newly inserted rows with the new product entries.
# Codes are entered in this particular order so the shared-string table grows
# in the same sequence as in the originally authored workbook
# (WB001, WQ001, WT001, WT002, WW001, WB002).
$ws.Range('B3785').Value = 'WB001'
$ws.Range('B3787').Value = 'WQ001'
$ws.Range('B3788').Value = 'WT001'
$ws.Range('B3789').Value = 'WT002'
$ws.Range('B3790').Value = 'WW001'
$ws.Range('B3786').Value = 'WB002'

$ws.Range('E3785').Value = 'Agua Oxigenada'
$ws.Range('K3785').Value = 2

$ws.Range('E3786').Value = 'Agua Oxigenada'
$ws.Range('K3786').Value = 2.5

$ws.Range('E3787').Value = 'Guante'
$ws.Range('K3787').Value = 12

$ws.Range('E3788').Value = 'Dispensador Ayudin'
$ws.Range('K3788').Value = 5

$ws.Range('E3789').Value = 'Esquinero Multiple'
$ws.Range('K3789').Value = 30

$ws.Range('E3790').Value = 'Cesto'
$ws.Range('K3790').Value = 15

# Leftover hidden filter-database name pointing at the Codigo column, a
# common artifact left behind in the workbook defined names collection.
$fd = $ws.Names.Add('_xlnm._FilterDatabase', '=Productos!$B$3:$B$3961')
$fd.Visible = $false

# Reset the view: drop the scrolled-down top-left cell and move the
# active selection to E3.
$ws.Range('E3').Select()
